# Applies the scheduled market-data refresh for the Belias_Profits sheets.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) for the
# affected leve rows on each job sheet, matching the latest Universalis pull.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H6").Value = 300.5
$ws.Range("I6").Value = 250
$ws.Range("J6").Value = 351
$ws.Range("K6").Value = 750
$ws.Range("L6").Value = 1053
$ws.Range("M6").Value = -638
$ws.Range("N6").Value = -1277

$ws.Range("H28").Value = 1369.0454
$ws.Range("I28").Value = 878.82355
$ws.Range("J28").Value = 3035.8
$ws.Range("K28").Value = 878.82355
$ws.Range("L28").Value = 3035.8
$ws.Range("M28").Value = -393.82355
$ws.Range("N28").Value = -4005.8

$ws.Range("H32").Value = 1115.6666
$ws.Range("J32").Value = 1139.8
$ws.Range("L32").Value = 1139.8
$ws.Range("N32").Value = -1791.8

$ws.Range("H41").Value = 318.05884
$ws.Range("I41").Value = 341.54544
$ws.Range("J41").Value = 275
$ws.Range("K41").Value = 341.54544
$ws.Range("L41").Value = 275
$ws.Range("M41").Value = 98.45456000000001
$ws.Range("N41").Value = -1155

$ws.Range("H62").Value = 4998.5625
$ws.Range("I62").Value = 3748.125
$ws.Range("J62").Value = 6249
$ws.Range("K62").Value = 3748.125
$ws.Range("L62").Value = 6249
$ws.Range("M62").Value = -3124.125
$ws.Range("N62").Value = -7497

$ws.Range("H64").Value = 3807.5
$ws.Range("I64").Value = 3563.3333
$ws.Range("J64").Value = 3888.889
$ws.Range("K64").Value = 3563.3333
$ws.Range("L64").Value = 3888.889
$ws.Range("M64").Value = -3315.3333
$ws.Range("N64").Value = -4384.889

$ws.Range("H65").Value = 4998.5625
$ws.Range("I65").Value = 3748.125
$ws.Range("J65").Value = 6249
$ws.Range("K65").Value = 18740.625
$ws.Range("L65").Value = 31245
$ws.Range("M65").Value = -15620.625
$ws.Range("N65").Value = -37485

$ws.Range("H67").Value = 3807.5
$ws.Range("I67").Value = 3563.3333
$ws.Range("J67").Value = 3888.889
$ws.Range("K67").Value = 3563.3333
$ws.Range("L67").Value = 3888.889
$ws.Range("M67").Value = -2705.3333
$ws.Range("N67").Value = -5604.889

$ws.Range("H76").Value = 3289.4736
$ws.Range("I76").Value = 3281.25
$ws.Range("J76").Value = 3333.3333
$ws.Range("K76").Value = 3281.25
$ws.Range("L76").Value = 3333.3333
$ws.Range("M76").Value = -2966.25
$ws.Range("N76").Value = -3963.3333

$ws.Range("H79").Value = 3289.4736
$ws.Range("I79").Value = 3281.25
$ws.Range("J79").Value = 3333.3333
$ws.Range("K79").Value = 3281.25
$ws.Range("L79").Value = 3333.3333
$ws.Range("M79").Value = -2189.25
$ws.Range("N79").Value = -5517.3333

$ws.Range("H98").Value = 1658.8334
$ws.Range("I98").Value = 1139.5385
$ws.Range("J98").Value = 2272.5454
$ws.Range("K98").Value = 1139.5385
$ws.Range("L98").Value = 2272.5454
$ws.Range("M98").Value = 358.4614999999999
$ws.Range("N98").Value = -5268.5454

$ws.Range("H122").Value = 1658.8334
$ws.Range("I122").Value = 1139.5385
$ws.Range("J122").Value = 2272.5454
$ws.Range("K122").Value = 3418.6155
$ws.Range("L122").Value = 6817.6362
$ws.Range("M122").Value = -968.6155000000003
$ws.Range("N122").Value = -11717.6362

$ws.Range("H132").Value = 16585497
$ws.Range("I132").Value = 2058453.8
$ws.Range("J132").Value = 76928600
$ws.Range("K132").Value = 6175361.4
$ws.Range("L132").Value = 230785800
$ws.Range("M132").Value = -6172831.4
$ws.Range("N132").Value = -230790860

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 4432.1064
$ws.Range("I32").Value = 2985.2307
$ws.Range("J32").Value = 13836.8
$ws.Range("K32").Value = 2985.2307
$ws.Range("L32").Value = 13836.8
$ws.Range("M32").Value = -2698.2307
$ws.Range("N32").Value = -14410.8

$ws.Range("H39").Value = 503
$ws.Range("I39").Value = 503
$ws.Range("K39").Value = 503
$ws.Range("M39").Value = 17

$ws.Range("H63").Value = 3575.75
$ws.Range("I63").Value = 2666.6667
$ws.Range("K63").Value = 2666.6667
$ws.Range("M63").Value = -1980.6667

$ws.Range("H66").Value = 3575.75
$ws.Range("I66").Value = 2666.6667
$ws.Range("K66").Value = 13333.3335
$ws.Range("M66").Value = -9901.333500000001

$ws.Range("H110").Value = 1205.3478
$ws.Range("I110").Value = 1236.1904
$ws.Range("J110").Value = 881.5
$ws.Range("K110").Value = 1236.1904
$ws.Range("L110").Value = 881.5
$ws.Range("M110").Value = 808.8096
$ws.Range("N110").Value = -4971.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 20001544
$ws.Range("I31").Value = 27027948
$ws.Range("K31").Value = 27027948
$ws.Range("M31").Value = -27027653

$ws.Range("H34").Value = 20001544
$ws.Range("I34").Value = 27027948
$ws.Range("K34").Value = 27027948
$ws.Range("M34").Value = -27027746

$ws.Range("H58").Value = 955.55554
$ws.Range("I58").Value = 855.2069
$ws.Range("J58").Value = 1371.2858
$ws.Range("K58").Value = 855.2069
$ws.Range("L58").Value = 1371.2858
$ws.Range("M58").Value = -652.2069
$ws.Range("N58").Value = -1777.2858

$ws.Range("H99").Value = 3088.1794
$ws.Range("I99").Value = 2962
$ws.Range("J99").Value = 3340.5386
$ws.Range("K99").Value = 2962
$ws.Range("L99").Value = 3340.5386
$ws.Range("M99").Value = -1464
$ws.Range("N99").Value = -6336.5386

$ws.Range("H126").Value = 3088.1794
$ws.Range("I126").Value = 2962
$ws.Range("J126").Value = 3340.5386
$ws.Range("K126").Value = 8886
$ws.Range("L126").Value = 10021.6158
$ws.Range("M126").Value = -6416
$ws.Range("N126").Value = -14961.6158

$ws.Range("H136").Value = 955.55554
$ws.Range("I136").Value = 855.2069
$ws.Range("J136").Value = 1371.2858
$ws.Range("K136").Value = 2565.6207
$ws.Range("L136").Value = 4113.857400000001
$ws.Range("M136").Value = -15.62069999999994
$ws.Range("N136").Value = -9213.857400000001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H109").Value = 2946.1538
$ws.Range("J109").Value = 3100
$ws.Range("L109").Value = 9300
$ws.Range("N109").Value = -11380

$ws.Range("H115").Value = 2371.182
$ws.Range("I115").Value = 1020.75
$ws.Range("J115").Value = 3142.8572
$ws.Range("K115").Value = 3062.25
$ws.Range("L115").Value = 9428.571599999999
$ws.Range("M115").Value = -1887.25
$ws.Range("N115").Value = -11778.5716

$ws.Range("H131").Value = 887.9
$ws.Range("J131").Value = 908.6129
$ws.Range("L131").Value = 2725.8387
$ws.Range("N131").Value = -12805.8387

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H113").Value = 1127.2667
$ws.Range("I113").Value = 968.3
$ws.Range("K113").Value = 968.3
$ws.Range("M113").Value = 1201.7

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H7").Value = 1500
$ws.Range("I7").Value = 1500
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1500
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1388
$ws.Range("N7").ClearContents()

$ws.Range("H126").Value = 1500
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2030
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 3444.0952
$ws.Range("I132").Value = 3429.7878
$ws.Range("J132").Value = 3496.5557
$ws.Range("K132").Value = 10289.3634
$ws.Range("L132").Value = 10489.6671
$ws.Range("M132").Value = -7759.3634
$ws.Range("N132").Value = -15549.6671

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H132").Value = 4084178.5
$ws.Range("I132").Value = 4446827.5
$ws.Range("K132").Value = 13340482.5
$ws.Range("M132").Value = -13337952.5
